$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Text columns (Date, Time, Weekday, Week) — leading "'" forces text so
# date/time-looking and numeric-looking strings are not auto-converted into
# a date serial / number, matching the existing rows' inlineStr text cells.
$ws.Cells.Item($row, 1).Value = "'2023-06-03"
$ws.Cells.Item($row, 2).Value = "'21:09:25"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "'22"

# Drop the "quote prefix" formatting flag the apostrophe trick leaves behind
# so the cells keep plain (unstyled) text, same as the rest of the sheet.
$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4)).ClearFormats()

# Numeric columns (Beijing .. Wuhan)
$ws.Cells.Item($row, 5).Value = 120881
$ws.Cells.Item($row, 6).Value = 133874
$ws.Cells.Item($row, 7).Value = 159090
$ws.Cells.Item($row, 8).Value = 130372
$ws.Cells.Item($row, 9).Value = 174723
$ws.Cells.Item($row, 10).Value = 112763
$ws.Cells.Item($row, 11).Value = 199825
$ws.Cells.Item($row, 12).Value = 218857
$ws.Cells.Item($row, 13).Value = 171846
$ws.Cells.Item($row, 14).Value = 118924
$ws.Cells.Item($row, 15).Value = 38197
$ws.Cells.Item($row, 16).Value = 34802
$ws.Cells.Item($row, 17).Value = 50250
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36976
$ws.Cells.Item($row, 20).Value = -1
